$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # Footer 1: Pearson logo, docPr id="3" -> rename image1.png to image2.png
    $f1 = $sec.Footers.Item(1)
    foreach ($shp in $f1.Range.InlineShapes) {
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image2.png"
        }
    }

    # Footer 2: Pearson logo, docPr id="2" -> rename image1.png to image2.png
    $f2 = $sec.Footers.Item(2)
    foreach ($shp in $f2.Range.InlineShapes) {
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image2.png"
        }
    }

    # Header 1 (second header, index 2 as wdHeaderFooterPrimary is blank): BTec logo -> rename image2.jpg to image1.jpg
    for ($i = 1; $i -le 3; $i++) {
        $hf = $sec.Headers.Item($i)
        if ($hf.Exists) {
            foreach ($shp in $hf.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }
}
